$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.29920000000001
$ws.Range("D4").Value = -8.582399999999998

$ws.Range("D5").Value = -8.577999999999998

$ws.Range("C6").Value = -11.6755
$ws.Range("D6").Value = -8.210099999999999

$ws.Range("C7").Value = -12.006

$ws.Range("C8").Value = -12.0604
$ws.Range("D8").Value = -8.046499999999998

$ws.Range("C16").Value = -11.6387
$ws.Range("D16").Value = -8.076200000000005

$ws.Range("C20").Value = -14.4487

$ws.Range("C21").Value = -12.97260000000001

$ws.Range("D22").Value = -7.964499999999997
